$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update in the order that reproduces the shared-string table ordering
$ws.Range("A2").Value = "basanagoudpatilas4573"
$ws.Range("A3").Value = "ananthreddithadi563546"
$ws.Range("B3").Value = "pass2@1234"
$ws.Range("C3").Value = "pass2@1234"
$ws.Range("B2").Value = "pass@1234"
$ws.Range("C2").Value = "pass@1234"
$ws.Range("F2").Value = "anandpatil@gmail.com"
$ws.Range("F3").Value = "Basupatil@gmail.com"

# Update selection to F3
$ws.Range("F3").Select()
